# Adds a "URL" column (C) to the dorm list sheet with a short slug per
# school, matching each row's existing Name/Nicknames entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
    "MIT","Stanford","Harvard","Caltech","UChicago","Princeton","Cornell","Yale",
    "Columbia","UPenn","Michigan","JHU","Duke","UC-Berkeley","UCLA","Northwestern",
    "UCSD","NYU","CMU","UWMadison","Brown","UT-Austin","UW","GIT","UIUC","UNC",
    "Rice","OSU","BU","PSU","WUSTL","Purdue","UC-Davis","USC","UMD","UCSB","Pitt",
    "MSU","Emory","UMN","UC-Irvine","UF","Dartmouth","Rochester","CWRU","Colorado",
    "UVA","Vanderbilt","TAMU","ASU","ND","UIC","Georgetown","Tufts","Miami",
    "Arizona","UMass","NCSU","Rutgers","Hawaii","YU","Buffalo","IUB","Northeastern",
    "UCSC","VirginiaTech","GWU","RPI","Utah","StonyBrook","Kansas","UConn","BC",
    "UCR","WakeForest","WSU","UCD","Tulane","UT-Knoxville","IIT","UT-Dallas",
    "Brandeis","UGA","Iowa","UDel","Wayne","CSU","OregonState","UMBC","Clark",
    "FSU","ISU","UNM","OU","Drexel","Lehigh","Howard","USF","NewSchool","UVM"
)

$ws.Range("C1").Value = "URL"
for ($i = 0; $i -lt $urls.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $urls[$i]
}

$ws.Columns.Item(3).ColumnWidth = 17.88671875

$ws.Range("C102").Select()
